# Update "new Madigan bike hours" - Riders and Average columns for the
# Ridership weekly sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Column C - Riders
$ws.Range("C2").Value = 197
$ws.Range("C4").Value = 206
$ws.Range("C5").Value = 160
$ws.Range("C6").Value = 229
$ws.Range("C8").Value = 80

# Column D - Average
$ws.Range("D2").Value = 227
$ws.Range("D3").Value = 218.23
$ws.Range("D4").Value = 211.81
$ws.Range("D5").Value = 236.6
$ws.Range("D6").Value = 241.77
$ws.Range("D7").Value = 113.15
$ws.Range("D8").Value = 91.95999999999999

$wb.Save()
